$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Задача 1"

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Человек заключил договор ренты на 20 лет, с ежемесячным платежом 30 тыс.руб. Первый платеж ожидается через месяц. Если стоимость денег равна 10% годовых, то сколько стоит такой договор? Сколько нужно заплатить в момент его заключения?"

$ws.Range("B1").Value = "n`nколичество периодов"
$ws.Range("C1").Value = "r`nстоимость денег"
$ws.Range("E1").Value = "PV`nцена договора в день заключения"
$ws.Range("F1").Value = "FV`nбудующая цена договора"
$ws.Range("D1").Value = "A`nмесячный платеж"

# Bold the short labels at the start of each header cell (B1,C1,D1,E1,F1)
$ws.Range("B1").Characters(1, 1).Font.Bold = $true
$ws.Range("C1").Characters(1, 1).Font.Bold = $true
$ws.Range("D1").Characters(1, 1).Font.Bold = $true
$ws.Range("E1").Characters(1, 2).Font.Bold = $true
$ws.Range("F1").Characters(1, 2).Font.Bold = $true

# Wrap text for headers
$ws.Range("A1:F1").WrapText = $true
$ws.Range("B1:F1").HorizontalAlignment = -4108  # xlCenter
$ws.Rows(1).RowHeight = 75

# ---- Data row (row 2) ----
$ws.Range("B2").Value = 20
$ws.Range("B2").HorizontalAlignment = -4108  # xlCenter

$ws.Range("C2").Value = 0.1
$ws.Range("C2").NumberFormat = "0%"

$ws.Range("D2").Value = 30000

$ws.Range("E2").Formula = "=PV(C2,B2,D2*12)"
$ws.Range("F2").Formula = "=FV(C2,B2,D2*12)"
$ws.Range("E2:F2").NumberFormat = "#,##0.00 ""₽"";[Red]-#,##0.00 ""₽"""

# ---- Column widths ----
$ws.Columns("A").ColumnWidth = 82.71
$ws.Columns("B").ColumnWidth = 11.86
$ws.Columns("C").ColumnWidth = 10.71
$ws.Columns("D").ColumnWidth = 10.43
$ws.Columns("E").ColumnWidth = 13.86
$ws.Columns("F").ColumnWidth = 14.86

# ---- Page setup ----
$ws.PageSetup.Orientation = 1  # xlPortrait
$ws.PageSetup.PaperSize = 9    # xlPaperA4
